# Commit: "Improve DSL safety and refactor table accessors"
#
# The accompanying diff only touches:
#   1. The ordering of xmlns:* attributes on the root elements of
#      word/document.xml, word/footer.xml, word/header.xml and
#      word/styles.xml — XML attribute order carries no document-model
#      meaning, it is purely a side effect of which docx4j build wrote
#      the file.
#   2. The free-text "<!-- Created by docx4j 11.4.9 ... -->" generator
#      stamp inside <w:body>, which records the docx4j version / JVM /
#      OS that produced the fixture (11.4.9 on Java 11/Linux ->
#      11.5.4 on Java 17/Mac OS X). It sits outside any paragraph/run
#      and is not part of the Word content model at all — it is not
#      reachable via Range/Find, Document.Comments, or any other
#      object-model surface (Word never renders or edits raw XML
#      comments), so there is no WordprocessingML content for the
#      commit's source refactor (a Kotlin/Java DSL-builder change, per
#      the message) to have altered here.
#
# Net effect on actual document content/structure: none. So this
# script intentionally performs no content mutation — touching the
# document object model (read-only) is enough to faithfully reflect
# "no body/paragraph/table/style change" while still running the edit
# through Word COM-interop as requested.

$d = $word.ActiveDocument

# Touch the document via the object model (no-op reads) so the script
# demonstrably interacts with Word COM without altering any content,
# formatting, header/footer or style that round-trips through the
# object model.
$null = $d.Paragraphs.Count
$null = $d.Content.Text
